$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (existing "pending" sale) -------------------------------------
# The sale that was previously the newest entry now gets the same explicit
# cell styling used by the other recorded rows (row 2-6), since it is no
# longer the "just typed" row once a newer sale (row 8) is recorded.
$ws.Range("A7:H7").Style = "Normal"

# --- Row 8 (new sale just processed) --------------------------------------
# Plain text / id fields - safe to assign directly.
$ws.Range("A8").Value = "95d64029-3b57-46ea-b0a6-8f8f87821cdc"
$ws.Range("B8").Value = "xksksjd"
$ws.Range("C8").Value = "Soft Drinks"
$ws.Range("D8").Value = "7UP"

# Quantity is a real number.
$ws.Range("E8").Value = 1

# Price/Date must be recorded as plain text (matching how every other row in
# this sheet stores these columns), not auto-converted to a number/date by
# Excel's smart entry. Force text via a scratch cell (NumberFormat "@") and
# bring the value across with a values-only paste so the destination cell's
# own formatting / the sheet's used range stay untouched. H8 is used as the
# scratch slot since it is about to receive its own real value anyway (it is
# already inside the new used range, so nothing extra gets marked as used).
$ws.Range("H8").NumberFormat = "@"

$ws.Range("H8").Value = "100"
$ws.Range("H8").Copy()
$ws.Range("F8").PasteSpecial(-4163)

$ws.Range("H8").Value = "2024-09-14"
$ws.Range("H8").Copy()
$ws.Range("G8").PasteSpecial(-4163)

# Reset the scratch cell back to an untouched state, then give it its real
# value - Time strings like "15:47:20" are not reinterpreted as numbers, so
# a plain assignment is enough here.
$ws.Range("H8").Clear()
$ws.Range("H8").Value = "15:47:20"
